$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.691.32"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").Value = "2.525.98"
$ws.Range("E3").Value = "  -2.42%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.567"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.52%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  -1.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.79%  "
$ws.Range("E11").Value = "  -1.20%  "
$ws.Range("E12").Value = "  -2.64%  "
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("D14").Value = "2.913.55"
$ws.Range("E14").Value = "  -2.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.67"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.70%  "
$ws.Range("D16").Value = "2.557.25"
$ws.Range("E16").Value = "  -0.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.807"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.66%  "
$ws.Range("D18").Value = "42.692.04"
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.32%  "
$ws.Range("E20").Value = "  -1.80%  "
$ws.Range("E21").Value = "  -4.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "243.77"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.07%  "
$ws.Range("E24").Value = "  -2.39%  "
$ws.Range("E25").Value = "  -2.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.87%  "
$ws.Range("E30").Value = "  -1.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "156.97"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.79"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.20%  "
$ws.Range("E33").Value = "  +10.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0787"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.40%  "
$ws.Range("E35").Value = "  -2.91%  "
$ws.Range("E36").Value = "  -5.89%  "
$ws.Range("E37").Value = "  -7.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.20"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.20%  "
$ws.Range("E39").Value = "  -1.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.119"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.20"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.09"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.59%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("E44").Value = "  -1.91%  "
$ws.Range("E45").Value = "  +1.46%  "
$ws.Range("D46").Value = "1.993.14"
$ws.Range("E46").Value = "  -1.10%  "
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("D48").Value = "2.766.93"
$ws.Range("E48").Value = "  -2.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "80.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.191"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.08%  "
